$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "
Category: CAT1  
Explanation: The incident involves multiple SSH brute force attempts from an unauthorized IP address, which directly relates to credential phishing (search term: `"phishing`") and unauthorized access to user/admin accounts. The log entries indicate attempts to exploit weak credentials, aligning with the definition of account compromise."
$ws.Range("B2").Value = "{'Category': 'CAT1', 'Explanation': 'The incident involves multiple SSH brute force attempts from an unauthorized IP address, which directly relates to credential phishing (search term: `"phishing`") and unauthorized access to user/admin accounts. The log entries indicate attempts to exploit weak credentials, aligning with the definition of account compromise.'}"
$ws.Range("C2").Value = "CAT1"
$ws.Range("D2").Value = "The incident involves multiple SSH brute force attempts from an unauthorized IP address, which directly relates to credential phishing (search term: `"phishing`") and unauthorized access to user/admin accounts. The log entries indicate attempts to exploit weak credentials, aligning with the definition of account compromise."

$ws.Range("A3").Value = "
Category: CAT1  
Explanation: The incident involves unauthorized access through SSH brute force, which is a classic example of Account Compromise. The log entries indicate multiple attempts to log in using invalid credentials, indicating a breach of user or administrator accounts."
$ws.Range("B3").Value = "{'Category': 'CAT1', 'Explanation': 'The incident involves unauthorized access through SSH brute force, which is a classic example of Account Compromise. The log entries indicate multiple attempts to log in using invalid credentials, indicating a breach of user or administrator accounts.'}"
$ws.Range("C3").Value = "CAT1"
$ws.Range("D3").Value = "The incident involves unauthorized access through SSH brute force, which is a classic example of Account Compromise. The log entries indicate multiple attempts to log in using invalid credentials, indicating a breach of user or administrator accounts."

$ws.Range("A4").Value = "
Category: CAT10  
Explanation: The incident involves misuse of resources by exploiting the IP address for non-authorized purposes, such as using the network for cryptocurrency mining or spam campaigns."
$ws.Range("B4").Value = "{'Category': 'CAT10', 'Explanation': 'The incident involves misuse of resources by exploiting the IP address for non-authorized purposes, such as using the network for cryptocurrency mining or spam campaigns.'}"
$ws.Range("C4").Value = "CAT10"
$ws.Range("D4").Value = "The incident involves misuse of resources by exploiting the IP address for non-authorized purposes, such as using the network for cryptocurrency mining or spam campaigns."

$ws.Range("A5").Value = "
Category: CAT1  
Explanation: The incident involves unauthorized access through SSH brute force attacks, which fall under Account Compromise as described in NIST category CAT1."
$ws.Range("B5").Value = "{'Category': 'CAT1', 'Explanation': 'The incident involves unauthorized access through SSH brute force attacks, which fall under Account Compromise as described in NIST category CAT1.'}"
$ws.Range("C5").Value = "CAT1"
$ws.Range("D5").Value = "The incident involves unauthorized access through SSH brute force attacks, which fall under Account Compromise as described in NIST category CAT1."

$ws.Range("A6").Value = "
Category: CAT4  
Explanation: The incident involves unauthorized disclosure of sensitive data through the SystemBC malware, which was used for lateral movement and data exfiltration."
$ws.Range("B6").Value = "{'Category': 'CAT4', 'Explanation': 'The incident involves unauthorized disclosure of sensitive data through the SystemBC malware, which was used for lateral movement and data exfiltration.'}"
$ws.Range("C6").Value = "CAT4"
$ws.Range("D6").Value = "The incident involves unauthorized disclosure of sensitive data through the SystemBC malware, which was used for lateral movement and data exfiltration."

$ws.Range("A7").Value = "
Category: CAT2  
Explanation: The incident involves the use of the malware SystemBC, which is classified under CAT2: Malware as an infection or malicious code. The description highlights its role as a backdoor and its application in ransomware attacks, aligning with the definition of malware."
$ws.Range("B7").Value = "{'Category': 'CAT2', 'Explanation': 'The incident involves the use of the malware SystemBC, which is classified under CAT2: Malware as an infection or malicious code. The description highlights its role as a backdoor and its application in ransomware attacks, aligning with the definition of malware.'}"
$ws.Range("C7").Value = "CAT2"
$ws.Range("D7").Value = "The incident involves the use of the malware SystemBC, which is classified under CAT2: Malware as an infection or malicious code. The description highlights its role as a backdoor and its application in ransomware attacks, aligning with the definition of malware."

$ws.Range("A8").Value = "
Category: CAT3  
Explanation: The incident involves a denial of service attack (DDoS) as evidenced by blocked ports and network disruptions, which aligns with NIST's definition of a DDoS attack."
$ws.Range("B8").Value = "{'Category': 'CAT3', 'Explanation': `"The incident involves a denial of service attack (DDoS) as evidenced by blocked ports and network disruptions, which aligns with NIST's definition of a DDoS attack.`"}"
$ws.Range("C8").Value = "CAT3"
$ws.Range("D8").Value = "The incident involves a denial of service attack (DDoS) as evidenced by blocked ports and network disruptions, which aligns with NIST's definition of a DDoS attack."

$ws.Range("A9").Value = "
Category: CAT11  
Explanation: The incident involves a third-party issue related to an abuse from an IP address, as indicated in the email content and the reference to the Abuse Contact Database."
$ws.Range("B9").Value = "{'Category': 'CAT11', 'Explanation': 'The incident involves a third-party issue related to an abuse from an IP address, as indicated in the email content and the reference to the Abuse Contact Database.'}"
$ws.Range("C9").Value = "CAT11"
$ws.Range("D9").Value = "The incident involves a third-party issue related to an abuse from an IP address, as indicated in the email content and the reference to the Abuse Contact Database."

$ws.Range("A10").Value = "
Category: CAT4  
Explanation: The incident describes unauthorized disclosure of sensitive data (leaked credentials) through the compromised machine, aligning with NIST's definition of a data leak."
$ws.Range("B10").Value = "{'Category': 'CAT4', 'Explanation': `"The incident describes unauthorized disclosure of sensitive data (leaked credentials) through the compromised machine, aligning with NIST's definition of a data leak.`"}"
$ws.Range("C10").Value = "CAT4"
$ws.Range("D10").Value = "The incident describes unauthorized disclosure of sensitive data (leaked credentials) through the compromised machine, aligning with NIST's definition of a data leak."

$ws.Range("A11").Value = "
Category: CAT12  
Explanation: The incident involves an unconfirmed or prevented intrusion attempt, as described in the response to the CERT email, which references a botnet attack using the provided IP address."
$ws.Range("B11").Value = "{'Category': 'CAT12', 'Explanation': 'The incident involves an unconfirmed or prevented intrusion attempt, as described in the response to the CERT email, which references a botnet attack using the provided IP address.'}"
$ws.Range("C11").Value = "CAT12"
$ws.Range("D11").Value = "The incident involves an unconfirmed or prevented intrusion attempt, as described in the response to the CERT email, which references a botnet attack using the provided IP address."

$ws.Range("A12").Value = "
Category: CAT3  
Explanation: The incident involves a coordinated DDoS botnet attack, which constitutes a denial of service attack (CAT3) as it disrupts network services and causes significant packet loss."
$ws.Range("B12").Value = "{'Category': 'CAT3', 'Explanation': 'The incident involves a coordinated DDoS botnet attack, which constitutes a denial of service attack (CAT3) as it disrupts network services and causes significant packet loss.'}"
$ws.Range("C12").Value = "CAT3"
$ws.Range("D12").Value = "The incident involves a coordinated DDoS botnet attack, which constitutes a denial of service attack (CAT3) as it disrupts network services and causes significant packet loss."

$ws.Range("A13").Value = "
Category: CAT4  
Explanation: The incident involves unauthorized disclosure of sensitive data (IP addresses) related to a DDoS botnet attack, which aligns with NIST's definition of a data leak."
$ws.Range("B13").Value = "{'Category': 'CAT4', 'Explanation': `"The incident involves unauthorized disclosure of sensitive data (IP addresses) related to a DDoS botnet attack, which aligns with NIST's definition of a data leak.`"}"
$ws.Range("C13").Value = "CAT4"
$ws.Range("D13").Value = "The incident involves unauthorized disclosure of sensitive data (IP addresses) related to a DDoS botnet attack, which aligns with NIST's definition of a data leak."

$ws.Range("A14").Value = "
Category: CAT5  
Explanation: The incident describes a vulnerability in the Zimbra Collaboration Suite, which falls under Vulnerability Exploitation (CAT5). The vulnerability is related to a CVE and exploits the software's flaws, allowing attackers to execute arbitrary code without authentication, as indicated in the incident description."
$ws.Range("B14").Value = "{'Category': 'CAT5', 'Explanation': `"The incident describes a vulnerability in the Zimbra Collaboration Suite, which falls under Vulnerability Exploitation (CAT5). The vulnerability is related to a CVE and exploits the software's flaws, allowing attackers to execute arbitrary code without authentication, as indicated in the incident description.`"}"
$ws.Range("C14").Value = "CAT5"
$ws.Range("D14").Value = "The incident describes a vulnerability in the Zimbra Collaboration Suite, which falls under Vulnerability Exploitation (CAT5). The vulnerability is related to a CVE and exploits the software's flaws, allowing attackers to execute arbitrary code without authentication, as indicated in the incident description."

$ws.Range("A15").Value = "
Category: CAT4  
Explanation: The incident involves unauthorized disclosure of sensitive data (leaked credentials) through the exposure of a BGP service. The service being exposed could be a data leak if the data is exposed to the public, which aligns with the `"data leak`" category."
$ws.Range("B15").Value = "{'Category': 'CAT4', 'Explanation': 'The incident involves unauthorized disclosure of sensitive data (leaked credentials) through the exposure of a BGP service. The service being exposed could be a data leak if the data is exposed to the public, which aligns with the `"data leak`" category.'}"
$ws.Range("C15").Value = "CAT4"
$ws.Range("D15").Value = "The incident involves unauthorized disclosure of sensitive data (leaked credentials) through the exposure of a BGP service. The service being exposed could be a data leak if the data is exposed to the public, which aligns with the `"data leak`" category."

$ws.Range("A16").Value = "
Category: CAT2  
Explanation: The incident involves unauthorized use of an IP address for port scans and database queries, which are indicative of a malware attack, specifically targeting vulnerabilities such as SQL injection or exploitation of existing systems."
$ws.Range("B16").Value = "{'Category': 'CAT2', 'Explanation': 'The incident involves unauthorized use of an IP address for port scans and database queries, which are indicative of a malware attack, specifically targeting vulnerabilities such as SQL injection or exploitation of existing systems.'}"
$ws.Range("C16").Value = "CAT2"
$ws.Range("D16").Value = "The incident involves unauthorized use of an IP address for port scans and database queries, which are indicative of a malware attack, specifically targeting vulnerabilities such as SQL injection or exploitation of existing systems."

$ws.Range("A17").Value = "
Category: CAT12  
Explanation: The incident involves an attempt to prevent or prevent attacks, specifically UDP-based amplification attacks, which is an intrusion attempt. The service being used for DDoS attacks is an attack being attempted."
$ws.Range("B17").Value = "{'Category': 'CAT12', 'Explanation': 'The incident involves an attempt to prevent or prevent attacks, specifically UDP-based amplification attacks, which is an intrusion attempt. The service being used for DDoS attacks is an attack being attempted.'}"
$ws.Range("C17").Value = "CAT12"
$ws.Range("D17").Value = "The incident involves an attempt to prevent or prevent attacks, specifically UDP-based amplification attacks, which is an intrusion attempt. The service being used for DDoS attacks is an attack being attempted."

$ws.Range("A18").Value = "
Category: Unknown  
Explanation: The incident involves a technical issue related to server configuration, but none of the predefined NIST categories directly apply to this scenario. The problem is about misconfigurations in NTP servers, which fall under cybersecurity but are not categorized in the provided list."
$ws.Range("B18").Value = "{'Category': 'Unknown', 'Explanation': 'The incident involves a technical issue related to server configuration, but none of the predefined NIST categories directly apply to this scenario. The problem is about misconfigurations in NTP servers, which fall under cybersecurity but are not categorized in the provided list.'}"
$ws.Range("C18").Value = "Unknown"
$ws.Range("D18").Value = "The incident involves a technical issue related to server configuration, but none of the predefined NIST categories directly apply to this scenario. The problem is about misconfigurations in NTP servers, which fall under cybersecurity but are not categorized in the provided list."

$ws.Range("A19").Value = "
Category: CAT7  
Explanation: The incident involves a phishing scam, which falls under social engineering (CAT7) as it involves deception to gain access or data."
$ws.Range("B19").Value = "{'Category': 'CAT7', 'Explanation': 'The incident involves a phishing scam, which falls under social engineering (CAT7) as it involves deception to gain access or data.'}"
$ws.Range("C19").Value = "CAT7"
$ws.Range("D19").Value = "The incident involves a phishing scam, which falls under social engineering (CAT7) as it involves deception to gain access or data."

$ws.Range("A20").Value = "
Category: CAT7  
Explanation: The incident involves phishing attempts to extort users, which falls under social engineering."
$ws.Range("B20").Value = "{'Category': 'CAT7', 'Explanation': 'The incident involves phishing attempts to extort users, which falls under social engineering.'}"
$ws.Range("C20").Value = "CAT7"
$ws.Range("D20").Value = "The incident involves phishing attempts to extort users, which falls under social engineering."

$ws.Range("A21").Value = "
Category: CAT4  
Explanation: The incident involves unauthorized disclosure of sensitive data (website content alteration) as described."
$ws.Range("B21").Value = "{'Category': 'CAT4', 'Explanation': 'The incident involves unauthorized disclosure of sensitive data (website content alteration) as described.'}"
$ws.Range("C21").Value = "CAT4"
$ws.Range("D21").Value = "The incident involves unauthorized disclosure of sensitive data (website content alteration) as described."

$ws.Range("A22").Value = "
Category: CAT4  
Explanation: The incident involves unauthorized disclosure of sensitive website content, which aligns with the NIST category of `"Data Leak.`" The alteration of the website's content indicates exposure of confidential data."
$ws.Range("B22").Value = "{'Category': 'CAT4', 'Explanation': 'The incident involves unauthorized disclosure of sensitive website content, which aligns with the NIST category of `"Data Leak.`" The alteration of the website\'s content indicates exposure of confidential data.'}"
$ws.Range("C22").Value = "CAT4"
$ws.Range("D22").Value = "The incident involves unauthorized disclosure of sensitive website content, which aligns with the NIST category of `"Data Leak.`" The alteration of the website's content indicates exposure of confidential data."

$ws.Range("A23").Value = "
Category: CAT4  
Explanation: The incident involves unauthorized disclosure of sensitive website content, as the alteration of the website's data could be a breach of confidentiality."
$ws.Range("B23").Value = "{'Category': 'CAT4', 'Explanation': `"The incident involves unauthorized disclosure of sensitive website content, as the alteration of the website's data could be a breach of confidentiality.`"}"
$ws.Range("C23").Value = "CAT4"
$ws.Range("D23").Value = "The incident involves unauthorized disclosure of sensitive website content, as the alteration of the website's data could be a breach of confidentiality."

$ws.Range("A24").Value = "
Category: CAT4  
Explanation: The incident involves unauthorized disclosure of sensitive website content, which aligns with the NIST category for `"Data Leak`" (unauthorized exposure of sensitive data). The altered website data falls under `"exposed data`" and `"leaked credentials,`" indicating a breach of confidentiality."
$ws.Range("B24").Value = "{'Category': 'CAT4', 'Explanation': 'The incident involves unauthorized disclosure of sensitive website content, which aligns with the NIST category for `"Data Leak`" (unauthorized exposure of sensitive data). The altered website data falls under `"exposed data`" and `"leaked credentials,`" indicating a breach of confidentiality.'}"
$ws.Range("C24").Value = "CAT4"
$ws.Range("D24").Value = "The incident involves unauthorized disclosure of sensitive website content, which aligns with the NIST category for `"Data Leak`" (unauthorized exposure of sensitive data). The altered website data falls under `"exposed data`" and `"leaked credentials,`" indicating a breach of confidentiality."

$ws.Range("A25").Value = "
Category: CAT4  
Explanation: The incident involves unauthorized disclosure of sensitive data (leaked credentials) through the website's altered content."
$ws.Range("B25").Value = "{'Category': 'CAT4', 'Explanation': `"The incident involves unauthorized disclosure of sensitive data (leaked credentials) through the website's altered content.`"}"
$ws.Range("C25").Value = "CAT4"
$ws.Range("D25").Value = "The incident involves unauthorized disclosure of sensitive data (leaked credentials) through the website's altered content."
